# "added light display bad"
# - Adds a new "drum" instruction sheet (positioned after machine2, before guitar)
# - Adds a matching "drum" Object row to the init sheet
# - Retimes / redirects the Light4 instruction sequence
# - Leaves Light4 as the active/selected sheet, Light3's selection nudged

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) init sheet: register the new "drum" object (row 17)
# ---------------------------------------------------------------------------
$init = $wb.Worksheets.Item("init")
$init.Activate()

$init.Range("A17").Value = "Object"
$init.Range("B17").Value = "drum"
$init.Range("C17").Value = "drum.png"
$init.Range("D17").Value = 200
$init.Range("E17").Value = 5
$init.Range("F17").Value = 150
$init.Range("G17").Value = 0

$init.Range("H17").Select()

# ---------------------------------------------------------------------------
# 2) Light3 sheet: selection moves to D5 (no data change)
# ---------------------------------------------------------------------------
$light3 = $wb.Worksheets.Item("Light3")
$light3.Activate()
$light3.Range("D5").Select()

# ---------------------------------------------------------------------------
# 3) Light4 sheet: re-time the moves and change their direction values
# ---------------------------------------------------------------------------
$light4 = $wb.Worksheets.Item("Light4")
$light4.Activate()

$light4.Range("B2").Value = 1
$light4.Range("C2").Value = -100

$light4.Range("B3").Value = 2
$light4.Range("C3").Value = -160

$light4.Range("B4").Value = 3
$light4.Range("C4").Value = -100

$light4.Range("B5").Value = 4
$light4.Range("C5").Value = -160

$light4.Range("B6").Value = 5
$light4.Range("C6").Value = -100

$light4.Range("B7").Value = 6
$light4.Range("C7").Value = -160

$light4.Range("C8").Select()

# ---------------------------------------------------------------------------
# 4) Insert a brand-new "drum" worksheet right after "machine2" and before
#    "guitar", carrying its own Move To / Loop To instruction table.
# ---------------------------------------------------------------------------
$machine2 = $wb.Worksheets.Item("machine2")
$drum = $wb.Worksheets.Add($null, $machine2)
$drum.Name = "drum"

$drum.Range("A1").Value = "Instruction"
$drum.Range("B1").Value = "End Time"
$drum.Range("C1").Value = "Horizontal Position"
$drum.Range("D1").Value = "Vertical Position"
$drum.Range("E1").Value = "Loop To Index"

$drum.Range("A2").Value = "Move To"
$drum.Range("B2").Value = 5
$drum.Range("C2").Value = 200
$drum.Range("D2").Value = 50

$drum.Range("A3").Value = "Move To"
$drum.Range("B3").Value = 10
$drum.Range("C3").Value = 100
$drum.Range("D3").Value = 0

$drum.Range("A4").Value = "Loop To"
$drum.Range("B4").Value = 10
$drum.Range("E4").Value = 0

$drum.Range("B5").Select()

# ---------------------------------------------------------------------------
# Leave Light4 as the final active sheet/tab (matches activeTab in the diff)
# ---------------------------------------------------------------------------
$light4.Activate()
$light4.Range("C8").Select()
